$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update existing section-header text in place (keeps shared-string slot) ---
$ws.Range("B3").Value2 = "code&comment2code -> Code Complexity (AST changes)"
$ws.Range("B7").Value2 = "code2comment -> Comment Complexity (# tokens)"

# --- 2. Build the new "code2comment -> Code Complexity (AST changes)" block (rows 12-15) ---
# 2a. merge the new section-header row first (while empty), then paste formats from B7:D7
$ws.Range("B12:D12").Merge()
$ws.Range("B7:D7").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2b. paste formats for the 3 data rows from rows 8-10
$ws.Range("B8:D10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2c. row heights
for ($i = 12; $i -le 15; $i++) {
    $ws.Rows.Item($i).RowHeight = 21.25
}

# 2d. values
$ws.Range("B12").Value2 = "code2comment -> Code Complexity (AST changes)"

$ws.Range("B13").Value2 = "All"
$ws.Range("C13").Value2 = [double]"4.434089e-56"
$ws.Range("C13").NumberFormat = "0.0#############################################################E+00"
$ws.Range("D13").Value2 = -0.1941853

$ws.Range("B14").Value2 = "T5CR"
$ws.Range("C14").Value2 = 0.0002513172
$ws.Range("C14").NumberFormat = "0.0#########E+00"
$ws.Range("D14").Value2 = -0.232217

$ws.Range("B15").Value2 = "CommentFinder"
$ws.Range("C15").Value2 = 0.00738188
$ws.Range("C15").NumberFormat = "0.0#######E+00"
$ws.Range("D15").Value2 = -0.1646046

"done"
